$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 86; existing rows 86.. shift down by one.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with its data.
$ws.Cells.Item(86, 1).Value = 4
$ws.Cells.Item(86, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(86, 3).Value = "Los Lagos"
$ws.Cells.Item(86, 4).Value = 44546
$ws.Cells.Item(86, 5).Value = 10
$ws.Cells.Item(86, 6).Value = 100112043
$ws.Cells.Item(86, 7).Value = "Pepino ensalada"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 160
$ws.Cells.Item(86, 11).Value = 10000
$ws.Cells.Item(86, 12).Value = 11000
$ws.Cells.Item(86, 13).Value = 10500
$ws.Cells.Item(86, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(86, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(86, 16).Value = 175
$ws.Cells.Item(86, 17).Value = 60
$ws.Cells.Item(86, 18).Value = "Hortaliza"
